$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for columns G (peakValue), H (RMS), I (tau), J (AUC) for rows 2-19.
$data = @(
    @(50.15000000000001,12.84536111111111,153.9534722222222,5027.915722222222),
    @(81.98425925925926,22.06077777777778,156.9860740740741,9333.02887037037),
    @(66.86018518518519,19.27598148148148,199.1854907407407,8785.400435185185),
    @(67.20555555555556,18.26519444444445,162.8410555555556,7508.073027777778),
    @(114.0638888888889,32.93819444444445,176.8305740740741,14621.79834259259),
    @(99.12037037037038,26.91130555555556,145.1042222222222,10822.44786111111),
    @(123.8351851851852,33.65185185185185,150.8492407407408,13712.78072222222),
    @(251.6518518518519,68.40933333333334,164.5025370370371,27733.00499074074),
    @(139.937962962963,33.95793518518519,97.69907407407408,11767.41068518518),
    @(143.3435185185185,38.61753703703704,137.7861018518519,14239.41121296296),
    @(216.7611111111111,49.80016666666667,69.07636111111111,15780.26272222222),
    @(184.6898148148148,50.48237037037037,161.3433888888889,20323.01669444445),
    @(105.9212962962963,26.42299074074074,116.863537037037,9238.739981481482),
    @(99.17129629629629,23.67728703703704,117.0621759259259,8837.139092592592),
    @(151.0166666666667,30.1268425925926,60.95931481481481,9303.403685185185),
    @(72.50555555555555,18.51413888888889,153.5498055555555,7226.979),
    @(122.9666666666667,34.91819444444445,144.6688333333333,14414.89663888889),
    @(119.6296296296296,29.74342592592593,121.2216944444444,10485.57519444444)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Range("G$row").Value = $vals[0]
    $ws.Range("H$row").Value = $vals[1]
    $ws.Range("I$row").Value = $vals[2]
    $ws.Range("J$row").Value = $vals[3]
}
